# Weekly fruit/vegetable price update: insert a new price record for
# "Granada" (Vega Modelo de Temuco) at row 69, pushing all subsequent
# records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 69 (shifts rows 69..164 down to 70..165).
$ws.Rows.Item(69).EntireRow.Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A69").Value = 10
$ws.Range("B69").Value = "Vega Modelo de Temuco"
$ws.Range("C69").Value = "La Araucanía"
$ws.Range("D69").Value = 44803
$ws.Range("E69").Value = 9
$ws.Range("F69").Value = "Fruta"
$ws.Range("G69").Value = 100104
$ws.Range("H69").Value = "Frutos de pepita"
$ws.Range("I69").Value = 100104001
$ws.Range("J69").Value = "Granada"
$ws.Range("K69").Value = "Wonderfull"
$ws.Range("L69").Value = "Primera"
$ws.Range("M69").Value = 160
$ws.Range("N69").Value = 14000
$ws.Range("O69").Value = 14000
$ws.Range("P69").Value = 14000
$ws.Range("Q69").Value = "$/bandeja 10 kilos granel"
$ws.Range("R69").Value = "Provincia de Limarí"
$ws.Range("S69").Value = 1400
$ws.Range("T69").Value = 10
